# Atualizacao de bases das ligas - corrige linhas de jogos trocadas
# (ids, equipas e respetivas odds/estatisticas) na folha "Poland IV Liga".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value2 = 6750018
$ws.Range("E6").Value2 = "Hutnik Warsaw"
$ws.Range("F6").Value2 = "Swit Starozreby"
$ws.Range("G6").Value2 = 1
$ws.Range("H6").Value2 = 1
$ws.Range("I6").Value2 = ""
$ws.Range("J6").Value2 = ""
$ws.Range("K6").Value2 = "D"
$ws.Range("L6").Value2 = 1.571
$ws.Range("N6").Value2 = 4.333
$ws.Range("O6").Value2 = 1.4
$ws.Range("P6").Value2 = 4.5
$ws.Range("Q6").Value2 = 6
$ws.Range("R6").Value2 = -1.25
$ws.Range("S6").Value2 = 1.8
$ws.Range("T6").Value2 = 2
$ws.Range("V6").Value2 = 1.8
$ws.Range("W6").Value2 = 2
$ws.Range("X6").Value2 = -1
$ws.Range("Y6").Value2 = 3.5
$ws.Range("AA6").Value2 = -1
$ws.Range("AB6").Value2 = 1
$ws.Range("AD6").Value2 = 1
$ws.Range("B7").Value2 = 6746871
$ws.Range("E7").Value2 = "WDA Swiecie"
$ws.Range("F7").Value2 = "Chemik Bydgoszcz"
$ws.Range("G7").Value2 = 2
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 2
$ws.Range("J7").Value2 = 0
$ws.Range("K7").Value2 = "H"
$ws.Range("L7").Value2 = 3.25
$ws.Range("N7").Value2 = 1.8
$ws.Range("O7").Value2 = 3.5
$ws.Range("P7").Value2 = 4
$ws.Range("Q7").Value2 = 1.727
$ws.Range("R7").Value2 = 0.75
$ws.Range("S7").Value2 = 1.825
$ws.Range("T7").Value2 = 1.975
$ws.Range("V7").Value2 = 1.9
$ws.Range("W7").Value2 = 1.9
$ws.Range("X7").Value2 = 2.5
$ws.Range("Y7").Value2 = -1
$ws.Range("AA7").Value2 = 0.825
$ws.Range("AB7").Value2 = -1
$ws.Range("AD7").Value2 = 0.8999999999999999
$ws.Range("F17").Value2 = "Korona Piaski"
$ws.Range("B30").Value2 = 7018460
$ws.Range("E30").Value2 = "Pomorzanin Torun"
$ws.Range("F30").Value2 = "Chemik Bydgoszcz"
$ws.Range("G30").Value2 = 0
$ws.Range("H30").Value2 = 5
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 2
$ws.Range("L30").Value2 = 3.9
$ws.Range("M30").Value2 = 4.5
$ws.Range("N30").Value2 = 1.571
$ws.Range("O30").Value2 = 5.5
$ws.Range("P30").Value2 = 5.5
$ws.Range("Q30").Value2 = 1.333
$ws.Range("R30").Value2 = 1.75
$ws.Range("S30").Value2 = 1.8
$ws.Range("T30").Value2 = 2
$ws.Range("U30").Value2 = 3.75
$ws.Range("V30").Value2 = 2
$ws.Range("W30").Value2 = 1.8
$ws.Range("Z30").Value2 = 0.333
$ws.Range("AA30").Value2 = -1
$ws.Range("AB30").Value2 = 1
$ws.Range("AC30").Value2 = 1
$ws.Range("B31").Value2 = 7021864
$ws.Range("E31").Value2 = "Marcovia Marki"
$ws.Range("F31").Value2 = "MKS Piaseczno"
$ws.Range("G31").Value2 = 2
$ws.Range("H31").Value2 = 3
$ws.Range("I31").Value2 = ""
$ws.Range("J31").Value2 = ""
$ws.Range("L31").Value2 = 2.5
$ws.Range("M31").Value2 = 3.4
$ws.Range("N31").Value2 = 2.375
$ws.Range("O31").Value2 = 4
$ws.Range("P31").Value2 = 3.8
$ws.Range("Q31").Value2 = 1.65
$ws.Range("R31").Value2 = 0.75
$ws.Range("S31").Value2 = 1.925
$ws.Range("T31").Value2 = 1.875
$ws.Range("U31").Value2 = 3
$ws.Range("V31").Value2 = 1.9
$ws.Range("W31").Value2 = 1.9
$ws.Range("Z31").Value2 = 0.6499999999999999
$ws.Range("AA31").Value2 = -0.5
$ws.Range("AB31").Value2 = 0.4375
$ws.Range("AC31").Value2 = 0.8999999999999999
$ws.Range("E32").Value2 = "WDA Swiecie"
$ws.Range("B41").Value2 = 7068598
$ws.Range("E41").Value2 = "Beskid Andrychow"
$ws.Range("F41").Value2 = "Niwa Nowa Wies"
$ws.Range("G41").Value2 = 3
$ws.Range("I41").Value2 = 1
$ws.Range("J41").Value2 = 1
$ws.Range("K41").Value2 = "H"
$ws.Range("L41").Value2 = 1.727
$ws.Range("N41").Value2 = 3.5
$ws.Range("O41").Value2 = 1.727
$ws.Range("Q41").Value2 = 3.5
$ws.Range("S41").Value2 = 1.775
$ws.Range("T41").Value2 = 2.025
$ws.Range("V41").Value2 = 1.975
$ws.Range("W41").Value2 = 1.825
$ws.Range("X41").Value2 = 0.7270000000000001
$ws.Range("Z41").Value2 = -1
$ws.Range("AA41").Value2 = 0.7749999999999999
$ws.Range("AB41").Value2 = -1
$ws.Range("AC41").Value2 = 0.9750000000000001
$ws.Range("AD41").Value2 = -1
$ws.Range("B42").Value2 = 7068601
$ws.Range("E42").Value2 = "Arka Gdynia II"
$ws.Range("F42").Value2 = "Anioly Garczegorze"
$ws.Range("G42").Value2 = 0
$ws.Range("I42").Value2 = 0
$ws.Range("J42").Value2 = 0
$ws.Range("K42").Value2 = "A"
$ws.Range("L42").Value2 = 1.95
$ws.Range("N42").Value2 = 2.875
$ws.Range("O42").Value2 = 1.85
$ws.Range("Q42").Value2 = 3.1
$ws.Range("S42").Value2 = 1.9
$ws.Range("T42").Value2 = 1.9
$ws.Range("V42").Value2 = 1.875
$ws.Range("W42").Value2 = 1.925
$ws.Range("X42").Value2 = -1
$ws.Range("Z42").Value2 = 2.1
$ws.Range("AA42").Value2 = -1
$ws.Range("AB42").Value2 = 0.8999999999999999
$ws.Range("AC42").Value2 = -1
$ws.Range("AD42").Value2 = 0.925
$ws.Range("B46").Value2 = 7068602
$ws.Range("E46").Value2 = "Tarnovia Tarnowo Podgorne"
$ws.Range("F46").Value2 = "Victoria Wrzesnia"
$ws.Range("J46").Value2 = 0
$ws.Range("L46").Value2 = 2.1
$ws.Range("M46").Value2 = 4.2
$ws.Range("N46").Value2 = 2.5
$ws.Range("O46").Value2 = 2.1
$ws.Range("P46").Value2 = 4.2
$ws.Range("Q46").Value2 = 2.55
$ws.Range("R46").Value2 = -0.25
$ws.Range("S46").Value2 = 1.95
$ws.Range("T46").Value2 = 1.85
$ws.Range("U46").Value2 = 3
$ws.Range("V46").Value2 = 1.775
$ws.Range("W46").Value2 = 2.025
$ws.Range("Y46").Value2 = 3.2
$ws.Range("AA46").Value2 = -0.5
$ws.Range("AB46").Value2 = 0.425
$ws.Range("AC46").Value2 = 0.7749999999999999
$ws.Range("B47").Value2 = 7068603
$ws.Range("E47").Value2 = "Korona Piaski"
$ws.Range("F47").Value2 = "Lipno Steszew"
$ws.Range("J47").Value2 = 1
$ws.Range("L47").Value2 = 3.75
$ws.Range("M47").Value2 = 4.333
$ws.Range("N47").Value2 = 1.615
$ws.Range("O47").Value2 = 3.75
$ws.Range("P47").Value2 = 4.333
$ws.Range("Q47").Value2 = 1.615
$ws.Range("R47").Value2 = 0.75
$ws.Range("S47").Value2 = 2
$ws.Range("T47").Value2 = 1.8
$ws.Range("U47").Value2 = 3.25
$ws.Range("V47").Value2 = 2.025
$ws.Range("W47").Value2 = 1.775
$ws.Range("Y47").Value2 = 3.333
$ws.Range("AA47").Value2 = 1
$ws.Range("AB47").Value2 = -1
$ws.Range("AC47").Value2 = 1.025
$ws.Range("F58").Value2 = "Marcovia Marki"
$ws.Range("E77").Value2 = "WDA Swiecie"
$ws.Range("E81").Value2 = "Pomorzanin Torun"
$ws.Range("F102").Value2 = "Korona Piaski"
$ws.Range("E108").Value2 = "Tarnovia Tarnowo Podgorne"
$ws.Range("E120").Value2 = "Marcovia Marki"
$ws.Range("F127").Value2 = "Pomorzanin Torun"
$ws.Range("B138").Value2 = 8061188
$ws.Range("E138").Value2 = "Termalica BB Nieciecza II"
$ws.Range("F138").Value2 = "Dalin Myslenice"
$ws.Range("G138").Value2 = 1
$ws.Range("H138").Value2 = 1
$ws.Range("I138").Value2 = 0
$ws.Range("K138").Value2 = "D"
$ws.Range("L138").Value2 = 1.909
$ws.Range("M138").Value2 = 3.75
$ws.Range("N138").Value2 = 3.1
$ws.Range("O138").Value2 = 1.909
$ws.Range("P138").Value2 = 3.75
$ws.Range("Q138").Value2 = 3.1
$ws.Range("R138").Value2 = -0.5
$ws.Range("U138").Value2 = 3.5
$ws.Range("X138").Value2 = -1
$ws.Range("Y138").Value2 = 2.75
$ws.Range("AA138").Value2 = -1
$ws.Range("AB138").Value2 = 0.825
$ws.Range("AC138").Value2 = -1
$ws.Range("AD138").Value2 = 0.9750000000000001
$ws.Range("B139").Value2 = 8061187
$ws.Range("E139").Value2 = "Korona Kielce II"
$ws.Range("F139").Value2 = "Spartakus Daleszyce"
$ws.Range("G139").Value2 = 5
$ws.Range("H139").Value2 = 0
$ws.Range("I139").Value2 = 3
$ws.Range("K139").Value2 = "H"
$ws.Range("L139").Value2 = 1.142
$ws.Range("M139").Value2 = 7
$ws.Range("N139").Value2 = 12
$ws.Range("O139").Value2 = 1.04
$ws.Range("P139").Value2 = 10
$ws.Range("Q139").Value2 = 29
$ws.Range("R139").Value2 = -3.5
$ws.Range("U139").Value2 = 4.25
$ws.Range("X139").Value2 = 0.04000000000000004
$ws.Range("Y139").Value2 = -1
$ws.Range("AA139").Value2 = 0.9750000000000001
$ws.Range("AB139").Value2 = -1
$ws.Range("AC139").Value2 = 0.825
$ws.Range("AD139").Value2 = -1
$ws.Range("E173").Value2 = "WDA Swiecie"
$ws.Range("F173").Value2 = "Pomorzanin Torun"
$ws.Range("F180").Value2 = "Korona Piaski"
$ws.Range("E191").Value2 = "WDA Swiecie"
$ws.Range("F197").Value2 = "Hutnik Warsaw"
$ws.Range("F205").Value2 = "WDA Swiecie"
